$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new data rows are appended after the last existing row (129).
# Re-use the formatting of row 129 (date style in column A, plain
# General style everywhere else) by copying it down first.
$ws.Range("A129:H129").Copy()
$ws.Range("A130:H130").PasteSpecial(-4122)
$ws.Range("A129:H129").Copy()
$ws.Range("A131:H131").PasteSpecial(-4122)

# --- Row 130 ---
$ws.Range("A130").Value = 45475.2916666667
$ws.Range("B130").Value = 0
$ws.Range("C130").Value = 2.5
$ws.Range("D130").Value = 2.5
$ws.Range("E130").Value = 2.5
$ws.Range("F130").Value = 2.5
# Column G holds numeric-looking text, so force text storage.
$ws.Range("G130").NumberFormat = "@"
$ws.Range("G130").Value = "2.5"
$ws.Range("H130").Value = "LS.MI"

# --- Row 131 ---
$ws.Range("A131").Value = 45476.6493402778
$ws.Range("B131").Value = 12000
$ws.Range("C131").Value = 2.5
$ws.Range("D131").Value = 2.30999994277954
$ws.Range("E131").Value = 2.5
$ws.Range("F131").Value = 2.30999994277954
$ws.Range("G131").NumberFormat = "@"
$ws.Range("G131").Value = "2.30999994277954"
$ws.Range("H131").Value = "LS.MI"

# Restore the plain (General) look of the text columns, now that the
# values are stored as text, so no stray number format is left visible
# on the cells themselves.
$ws.Range("G129:H129").Copy()
$ws.Range("G130:H130").PasteSpecial(-4122)
$ws.Range("G129:H129").Copy()
$ws.Range("G131:H131").PasteSpecial(-4122)
